$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Spring 24 week 15 inputs - updated matchup averages

$ws.Range("D3").Value = 10.35
$ws.Range("E3").Value = 10.75

$ws.Range("C4").Value = 9.65
$ws.Range("E4").Value = 10.69
$ws.Range("F4").Value = 10.12

$ws.Range("C5").Value = 9.25
$ws.Range("D5").Value = 9.31
$ws.Range("F5").Value = 10.19
$ws.Range("H5").Value = 8.699999999999999

$ws.Range("D6").Value = 9.880000000000001
$ws.Range("E6").Value = 9.81
$ws.Range("G6").Value = 10.3
$ws.Range("J6").Value = 7.6

$ws.Range("F7").Value = 9.699999999999999
$ws.Range("H7").Value = 9.869999999999999
$ws.Range("I7").Value = 5.91

$ws.Range("E8").Value = 11.3
$ws.Range("G8").Value = 10.13
$ws.Range("I8").Value = 9.08

$ws.Range("G9").Value = 14.09
$ws.Range("H9").Value = 10.92

$ws.Range("F10").Value = 12.4
